# Auto-generated edit script applying numeric updates across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N17").Value = -1315.5
$ws.Range("J17").Value = 326.5
$ws.Range("L17").Value = 979.5
$ws.Range("H17").Value = 326.5
$ws.Range("N19").Value = -2969.182
$ws.Range("K19").Value = 1845.2858
$ws.Range("I19").Value = 1845.2858
$ws.Range("M19").Value = -1670.2858
$ws.Range("J19").Value = 2619.182
$ws.Range("L19").Value = 2619.182
$ws.Range("H19").Value = 2185.8
$ws.Range("K21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H21").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("N26").Value = -25688
$ws.Range("J26").Value = 25000
$ws.Range("L26").Value = 25000
$ws.Range("H26").Value = 25000
$ws.Range("N28").Value = -4346.5
$ws.Range("K28").Value = 293.16666
$ws.Range("I28").Value = 293.16666
$ws.Range("M28").Value = 191.83334
$ws.Range("J28").Value = 3376.5
$ws.Range("L28").Value = 3376.5
$ws.Range("H28").Value = 1526.5
$ws.Range("N29").Value = -33262
$ws.Range("K29").Value = 13099.0005
$ws.Range("I29").Value = 4366.3335
$ws.Range("M29").Value = -12818.0005
$ws.Range("J29").Value = 10900
$ws.Range("L29").Value = 32700
$ws.Range("H29").Value = 5999.75
$ws.Range("K33").Value = 282.51514
$ws.Range("I33").Value = 282.51514
$ws.Range("M33").Value = -53.51513999999997
$ws.Range("H33").Value = 771.9737
$ws.Range("K38").Value = 228.42858
$ws.Range("I38").Value = 76.14286
$ws.Range("M38").Value = 143.57142
$ws.Range("H38").Value = 76.14286
$ws.Range("K53").Value = 164.25
$ws.Range("I53").Value = 164.25
$ws.Range("M53").Value = 472.75
$ws.Range("H53").Value = 131.16667
$ws.Range("N87").ClearContents()
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("N111").Value = -29639
$ws.Range("K111").Value = 10354.8
$ws.Range("I111").Value = 3451.6
$ws.Range("M111").Value = -7287.799999999999
$ws.Range("J111").Value = 7835
$ws.Range("L111").Value = 23505
$ws.Range("H111").Value = 4846.3184
$ws.Range("N113").Value = -11715.8887
$ws.Range("K113").Value = 9233.333000000001
$ws.Range("I113").Value = 9233.333000000001
$ws.Range("M113").Value = -5979.333000000001
$ws.Range("J113").Value = 5207.8887
$ws.Range("L113").Value = 5207.8887
$ws.Range("H113").Value = 6214.25
$ws.Range("K125").Value = 8871.75
$ws.Range("I125").Value = 985.75
$ws.Range("M125").Value = -6411.75
$ws.Range("H125").Value = 978.6
$ws.Range("N135").Value = -125625.999
$ws.Range("K135").Value = 9288.473399999999
$ws.Range("I135").Value = 1032.0526
$ws.Range("M135").Value = -6753.473399999999
$ws.Range("J135").Value = 13395.111
$ws.Range("L135").Value = 120555.999
$ws.Range("H135").Value = 5005.893
$ws.Range("N136").Value = -80195
$ws.Range("J136").Value = 69995
$ws.Range("L136").Value = 69995
$ws.Range("H136").Value = 69995
$ws.Range("N138").Value = -23570.162
$ws.Range("K138").Value = 4357.0002
$ws.Range("I138").Value = 1452.3334
$ws.Range("M138").Value = 782.9997999999996
$ws.Range("J138").Value = 4430.054
$ws.Range("L138").Value = 13290.162
$ws.Range("H138").Value = 3847.4565
$ws.Range("K141").Value = 12239.1
$ws.Range("I141").Value = 4079.7
$ws.Range("M141").Value = -7059.099999999999
$ws.Range("H141").Value = 3990.0833

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N19").ClearContents()
$ws.Range("K19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("H19").Value = 0
$ws.Range("N45").Value = -5586.5
$ws.Range("K45").Value = 13887
$ws.Range("I45").Value = 13887
$ws.Range("M45").Value = -13510
$ws.Range("J45").Value = 4832.5
$ws.Range("L45").Value = 4832.5
$ws.Range("H45").Value = 10533.481
$ws.Range("N74").Value = -3600
$ws.Range("K74").Value = 2006.5
$ws.Range("I74").Value = 2006.5
$ws.Range("M74").Value = -1132.5
$ws.Range("J74").Value = 1852
$ws.Range("L74").Value = 1852
$ws.Range("H74").Value = 1984.4286
$ws.Range("N77").Value = -17996
$ws.Range("K77").Value = 10032.5
$ws.Range("I77").Value = 2006.5
$ws.Range("M77").Value = -5664.5
$ws.Range("J77").Value = 1852
$ws.Range("L77").Value = 9260
$ws.Range("H77").Value = 1984.4286
$ws.Range("N101").Value = -66489.33199999999
$ws.Range("J101").Value = 59999.332
$ws.Range("L101").Value = 59999.332
$ws.Range("H101").Value = 59999.332
$ws.Range("K102").Value = 1686.1818
$ws.Range("I102").Value = 1686.1818
$ws.Range("M102").Value = -64.18180000000007
$ws.Range("H102").Value = 1659.4642
$ws.Range("K122").Value = 8923.5879
$ws.Range("I122").Value = 2974.5293
$ws.Range("M122").Value = -6473.5879
$ws.Range("H122").Value = 3198.9443
$ws.Range("N132").Value = -10706.4998
$ws.Range("K132").Value = 6947.499899999999
$ws.Range("I132").Value = 2315.8333
$ws.Range("M132").Value = -4417.499899999999
$ws.Range("J132").Value = 1882.1666
$ws.Range("L132").Value = 5646.4998
$ws.Range("H132").Value = 2207.4167
$ws.Range("N138").Value = -98687.836
$ws.Range("J138").Value = 88407.836
$ws.Range("L138").Value = 88407.836
$ws.Range("H138").Value = 88407.836
$ws.Range("N141").Value = -100359.664
$ws.Range("J141").Value = 89999.664
$ws.Range("L141").Value = 89999.664
$ws.Range("H141").Value = 89999.664

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N81").Value = -37121
$ws.Range("J81").Value = 34999
$ws.Range("L81").Value = 34999
$ws.Range("H81").Value = 34999
$ws.Range("N84").Value = -115605
$ws.Range("J84").Value = 34999
$ws.Range("L84").Value = 104997
$ws.Range("H84").Value = 34999
$ws.Range("N99").Value = -8245.5
$ws.Range("K99").Value = 4069.077
$ws.Range("I99").Value = 4069.077
$ws.Range("M99").Value = -2571.077
$ws.Range("J99").Value = 5249.5
$ws.Range("L99").Value = 5249.5
$ws.Range("H99").Value = 4226.467
$ws.Range("K105").Value = 2168.6924
$ws.Range("I105").Value = 2168.6924
$ws.Range("M105").Value = -421.6923999999999
$ws.Range("H105").Value = 2146.1333
$ws.Range("N107").Value = -7342.2
$ws.Range("K107").Value = 2529.1
$ws.Range("I107").Value = 2529.1
$ws.Range("M107").Value = -609.0999999999999
$ws.Range("J107").Value = 3502.2
$ws.Range("L107").Value = 3502.2
$ws.Range("H107").Value = 2853.4666
$ws.Range("N133").Value = -80117.75
$ws.Range("J133").Value = 69997.75
$ws.Range("L133").Value = 69997.75
$ws.Range("H133").Value = 69997.75
$ws.Range("N134").Value = -12194.25
$ws.Range("K134").Value = 6231.615
$ws.Range("I134").Value = 2077.205
$ws.Range("M134").Value = -3696.615
$ws.Range("J134").Value = 2374.75
$ws.Range("L134").Value = 7124.25
$ws.Range("H134").Value = 2104.8838

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N16").Value = -24532.285
$ws.Range("K16").Value = 1152.1666
$ws.Range("I16").Value = 1152.1666
$ws.Range("M16").Value = -865.1666
$ws.Range("J16").Value = 23958.285
$ws.Range("L16").Value = 23958.285
$ws.Range("H16").Value = 13432.385
$ws.Range("N43").Value = -67529.25
$ws.Range("K43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("J43").Value = 67161.25
$ws.Range("L43").Value = 67161.25
$ws.Range("H43").Value = 67161.25
$ws.Range("N58").ClearContents()
$ws.Range("K58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("N86").Value = -10621.667
$ws.Range("K86").Value = 30309504
$ws.Range("I86").Value = 30309504
$ws.Range("M86").Value = -30308381
$ws.Range("J86").Value = 8375.666999999999
$ws.Range("L86").Value = 8375.666999999999
$ws.Range("H86").Value = 16673996
$ws.Range("N89").Value = -53110.335
$ws.Range("K89").Value = 151547520
$ws.Range("I89").Value = 30309504
$ws.Range("M89").Value = -151541904
$ws.Range("J89").Value = 8375.666999999999
$ws.Range("L89").Value = 41878.335
$ws.Range("H89").Value = 16673996
$ws.Range("N101").Value = -73651.25
$ws.Range("K101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("J101").Value = 67161.25
$ws.Range("L101").Value = 67161.25
$ws.Range("H101").Value = 67161.25
$ws.Range("N113").Value = -28298.285
$ws.Range("K113").Value = 1152.1666
$ws.Range("I113").Value = 1152.1666
$ws.Range("M113").Value = 1017.8334
$ws.Range("J113").Value = 23958.285
$ws.Range("L113").Value = 23958.285
$ws.Range("H113").Value = 13432.385
$ws.Range("N123").Value = -99793.5
$ws.Range("J123").Value = 89993.5
$ws.Range("L123").Value = 89993.5
$ws.Range("H123").Value = 89993.5
$ws.Range("K134").Value = 12509.6001
$ws.Range("I134").Value = 4169.8667
$ws.Range("M134").Value = -9974.6001
$ws.Range("H134").Value = 4221.75
$ws.Range("N136").ClearContents()
$ws.Range("K136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("H136").Value = 0
$ws.Range("N138").Value = -86727.39999999999
$ws.Range("J138").Value = 76447.39999999999
$ws.Range("L138").Value = 76447.39999999999
$ws.Range("H138").Value = 76447.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N4").Value = -6524
$ws.Range("J4").Value = 2100
$ws.Range("L4").Value = 6300
$ws.Range("H4").Value = 58894308
$ws.Range("N12").Value = -1371.78568
$ws.Range("K12").Value = 1395.42855
$ws.Range("I12").Value = 465.14285
$ws.Range("M12").Value = -1222.42855
$ws.Range("J12").Value = 341.92856
$ws.Range("L12").Value = 1025.78568
$ws.Range("H12").Value = 383
$ws.Range("K51").Value = 7617.8568
$ws.Range("I51").Value = 2539.2856
$ws.Range("M51").Value = -7157.8568
$ws.Range("H51").Value = 2539.2856
$ws.Range("N113").Value = -6371.3333
$ws.Range("K113").Value = 2257.8
$ws.Range("I113").Value = 752.6
$ws.Range("M113").Value = -87.80000000000018
$ws.Range("J113").Value = 677.1111
$ws.Range("L113").Value = 2031.3333
$ws.Range("H113").Value = 716.8421
$ws.Range("N118").Value = -26483
$ws.Range("K118").Value = 9299.25
$ws.Range("I118").Value = 3099.75
$ws.Range("M118").Value = -8056.25
$ws.Range("J118").Value = 7999
$ws.Range("L118").Value = 23997
$ws.Range("H118").Value = 4079.6
$ws.Range("K121").Value = 1027.5
$ws.Range("I121").Value = 342.5
$ws.Range("M121").Value = 282.5
$ws.Range("H121").Value = 238723.95
$ws.Range("N133").Value = -68244.25
$ws.Range("K133").Value = 37491
$ws.Range("I133").Value = 12497
$ws.Range("M133").Value = -32431
$ws.Range("J133").Value = 19374.75
$ws.Range("L133").Value = 58124.25
$ws.Range("H133").Value = 17999.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("K48").Value = 6000
$ws.Range("I48").Value = 6000
$ws.Range("M48").Value = -5515
$ws.Range("H48").Value = 6000
$ws.Range("N80").Value = -19994
$ws.Range("K80").Value = 1802.6
$ws.Range("I80").Value = 1802.6
$ws.Range("M80").Value = -804.5999999999999
$ws.Range("J80").Value = 17998
$ws.Range("L80").Value = 17998
$ws.Range("H80").Value = 3274.9092
$ws.Range("N83").Value = -99974
$ws.Range("K83").Value = 9013
$ws.Range("I83").Value = 1802.6
$ws.Range("M83").Value = -4021
$ws.Range("J83").Value = 17998
$ws.Range("L83").Value = 89990
$ws.Range("H83").Value = 3274.9092
$ws.Range("N101").Value = -33436.5
$ws.Range("J101").Value = 26946.5
$ws.Range("L101").Value = 26946.5
$ws.Range("H101").Value = 26946.5
$ws.Range("N102").ClearContents()
$ws.Range("K102").Value = 2999
$ws.Range("I102").Value = 2999
$ws.Range("M102").Value = -1377
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("H102").Value = 2999
$ws.Range("N112").Value = -49716
$ws.Range("J112").Value = 47500
$ws.Range("L112").Value = 47500
$ws.Range("H112").Value = 47500
$ws.Range("K122").Value = 10705.8465
$ws.Range("I122").Value = 3568.6155
$ws.Range("M122").Value = -8255.8465
$ws.Range("H122").Value = 4198.6875
$ws.Range("N126").Value = -16910
$ws.Range("K126").Value = 11959.5
$ws.Range("I126").Value = 3986.5
$ws.Range("M126").Value = -9489.5
$ws.Range("J126").Value = 3990
$ws.Range("L126").Value = 11970
$ws.Range("H126").Value = 3987.6667
$ws.Range("N130").Value = -74938.66800000001
$ws.Range("J130").Value = 64898.668
$ws.Range("L130").Value = 64898.668
$ws.Range("H130").Value = 64898.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N13").ClearContents()
$ws.Range("K13").Value = 6000
$ws.Range("I13").Value = 6000
$ws.Range("M13").Value = -5860
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("H13").Value = 6000
$ws.Range("N111").Value = -59526
$ws.Range("J111").Value = 51346
$ws.Range("L111").Value = 51346
$ws.Range("H111").Value = 51346
$ws.Range("N132").Value = -13478.1875
$ws.Range("K132").Value = 8493.828
$ws.Range("I132").Value = 2831.276
$ws.Range("M132").Value = -5963.828
$ws.Range("J132").Value = 2806.0625
$ws.Range("L132").Value = 8418.1875
$ws.Range("H132").Value = 2825.8242
$ws.Range("N136").Value = -13197
$ws.Range("K136").Value = 3802.875
$ws.Range("I136").Value = 1267.625
$ws.Range("M136").Value = -1252.875
$ws.Range("J136").Value = 2699
$ws.Range("L136").Value = 8097
$ws.Range("H136").Value = 1426.6666
$ws.Range("N138").Value = -79279
$ws.Range("J138").Value = 68999
$ws.Range("L138").Value = 68999
$ws.Range("H138").Value = 68999
$ws.Range("N139").Value = -70995
$ws.Range("J139").Value = 60715
$ws.Range("L139").Value = 60715
$ws.Range("H139").Value = 60715
$ws.Range("N141").Value = -87817.39999999999
$ws.Range("J141").Value = 77457.39999999999
$ws.Range("L141").Value = 77457.39999999999
$ws.Range("H141").Value = 77457.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N2").Value = -7223
$ws.Range("K2").Value = 4999
$ws.Range("I2").Value = 4999
$ws.Range("M2").Value = -4887
$ws.Range("J2").Value = 6999
$ws.Range("L2").Value = 6999
$ws.Range("H2").Value = 5665.6665
$ws.Range("N16").Value = -50584
$ws.Range("J16").Value = 50000
$ws.Range("L16").Value = 50000
$ws.Range("H16").Value = 50000
$ws.Range("K21").Value = 50000
$ws.Range("I21").Value = 50000
$ws.Range("M21").Value = -49765
$ws.Range("H21").Value = 66666.336
$ws.Range("K35").Value = 50000
$ws.Range("I35").Value = 50000
$ws.Range("M35").Value = -49710
$ws.Range("H35").Value = 66666.336
$ws.Range("N45").Value = -20533.2
$ws.Range("K45").Value = 45975.75
$ws.Range("I45").Value = 45975.75
$ws.Range("M45").Value = -45484.75
$ws.Range("J45").Value = 19551.2
$ws.Range("L45").Value = 19551.2
$ws.Range("H45").Value = 31295.445
$ws.Range("K62").Value = 17498.5
$ws.Range("I62").Value = 17498.5
$ws.Range("M62").Value = -16874.5
$ws.Range("H62").Value = 17498.5
$ws.Range("K65").Value = 87492.5
$ws.Range("I65").Value = 17498.5
$ws.Range("M65").Value = -84372.5
$ws.Range("H65").Value = 17498.5
$ws.Range("N81").Value = -3702.5
$ws.Range("J81").Value = 790.25
$ws.Range("L81").Value = 1580.5
$ws.Range("H81").Value = 1835.5454
$ws.Range("N84").Value = -18510.5
$ws.Range("J84").Value = 790.25
$ws.Range("L84").Value = 7902.5
$ws.Range("H84").Value = 1835.5454
$ws.Range("N107").Value = -24582.75
$ws.Range("K107").Value = 12226.5
$ws.Range("I107").Value = 4075.5
$ws.Range("M107").Value = -10306.5
$ws.Range("J107").Value = 6914.25
$ws.Range("L107").Value = 20742.75
$ws.Range("H107").Value = 4785.1875
$ws.Range("N132").Value = -19526.3339
$ws.Range("K132").Value = 8297.222099999999
$ws.Range("I132").Value = 2765.7407
$ws.Range("M132").Value = -5767.222099999999
$ws.Range("J132").Value = 4822.1113
$ws.Range("L132").Value = 14466.3339
$ws.Range("H132").Value = 3279.8333
$ws.Range("N136").Value = -30285
$ws.Range("K136").Value = 4153.5
$ws.Range("I136").Value = 1384.5
$ws.Range("M136").Value = -1603.5
$ws.Range("J136").Value = 8395
$ws.Range("L136").Value = 25185
$ws.Range("H136").Value = 1626.2413
$ws.Range("N138").Value = -108938.336
$ws.Range("J138").Value = 98658.336
$ws.Range("L138").Value = 98658.336
$ws.Range("H138").Value = 98658.336

